$d = $word.ActiveDocument

# Locate the bullet paragraph "Statut doit être implanté comme une liste fixe"
# (it is split across multiple runs: "...co" + "m" + "me une liste fixe").
# Find.Execute matches across run boundaries, returning a Range over the hit.
$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("Statut doit être implanté comme une liste fixe", `
                              $true, $false, $false, $false, $false, $true, 1, `
                              $false, "", 0)

if ($found) {
    # Grab the whole paragraph (including its paragraph mark) that holds the
    # matched text and delete it outright, so the bullet list item disappears
    # entirely rather than leaving an empty list entry behind.
    $para = $range.Paragraphs(1)
    $para.Range.Delete()
}
